$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 2020
$ws.Range("I5").Value = 1287.4000000000001
$ws.Range("I6").Value = 56.6
$ws.Range("I8").Value = 2.5
$ws.Range("I9").Value = 9.3000000000000007
$ws.Range("I10").Value = 0.9
